{"js": "// Update the answers in the \"two-digit division\" practice table.\n//\n// The table cells hold strings shaped like \"A\u00f7B=Q, R\" (dividend\u00f7divisor=\n// quotient, remainder). Each entry below is [oldText, newText]: every\n// oldText value is a unique, exact string within the document body, so a\n// plain exact-match search/replace per pair is sufficient and the pairs\n// can be applied in any order without collisions.\nconst replacements = [\n  [\"54\u00f73=18, 0\", \"23\u00f72=11, 1\"],\n  [\"10\u00f74=2, 2\", \"64\u00f79=7, 1\"],\n  [\"18\u00f76=3, 0\", \"25\u00f79=2, 7\"],\n  [\"25\u00f73=8, 1\", \"88\u00f72=44, 0\"],\n  [\"14\u00f73=4, 2\", \"22\u00f73=7, 1\"],\n  [\"54\u00f74=13, 2\", \"40\u00f79=4, 4\"],\n  [\"16\u00f79=1, 7\", \"11\u00f76=1, 5\"],\n  [\"74\u00f72=37, 0\", \"47\u00f79=5, 2\"],\n  [\"64\u00f73=21, 1\", \"49\u00f72=24, 1\"],\n  [\"49\u00f78=6, 1\", \"14\u00f75=2, 4\"],\n  [\"50\u00f74=12, 2\", \"43\u00f78=5, 3\"],\n  [\"95\u00f76=15, 5\", \"28\u00f73=9, 1\"],\n  [\"31\u00f77=4, 3\", \"33\u00f76=5, 3\"],\n  [\"87\u00f79=9, 6\", \"95\u00f73=31, 2\"],\n  [\"41\u00f76=6, 5\", \"32\u00f78=4, 0\"],\n  [\"21\u00f73=7, 0\", \"12\u00f74=3, 0\"],\n  [\"89\u00f79=9, 8\", \"61\u00f75=12, 1\"],\n  [\"86\u00f78=10, 6\", \"90\u00f75=18, 0\"],\n  [\"99\u00f73=33, 0\", \"58\u00f77=8, 2\"],\n  [\"11\u00f75=2, 1\", \"93\u00f76=15, 3\"],\n  [\"33\u00f77=4, 5\", \"16\u00f72=8, 0\"],\n  [\"14\u00f76=2, 2\", \"34\u00f76=5, 4\"],\n  [\"44\u00f72=22, 0\", \"38\u00f76=6, 2\"],\n  [\"19\u00f76=3, 1\", \"72\u00f72=36, 0\"],\n  [\"22\u00f77=3, 1\", \"51\u00f78=6, 3\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace specific two-digit-division answer cells in the table.\n# Each pair is (oldText, newText); all oldText values are unique, exact\n# strings within the document, so Find/Replace with exact matching and\n# no wildcards is sufficient and order-independent.\n$pairs = @(\n    @(\"54\u00f73=18, 0\", \"23\u00f72=11, 1\"),\n    @(\"10\u00f74=2, 2\", \"64\u00f79=7, 1\"),\n    @(\"18\u00f76=3, 0\", \"25\u00f79=2, 7\"),\n    @(\"25\u00f73=8, 1\", \"88\u00f72=44, 0\"),\n    @(\"14\u00f73=4, 2\", \"22\u00f73=7, 1\"),\n    @(\"54\u00f74=13, 2\", \"40\u00f79=4, 4\"),\n    @(\"16\u00f79=1, 7\", \"11\u00f76=1, 5\"),\n    @(\"74\u00f72=37, 0\", \"47\u00f79=5, 2\"),\n    @(\"64\u00f73=21, 1\", \"49\u00f72=24, 1\"),\n    @(\"49\u00f78=6, 1\", \"14\u00f75=2, 4\"),\n    @(\"50\u00f74=12, 2\", \"43\u00f78=5, 3\"),\n    @(\"95\u00f76=15, 5\", \"28\u00f73=9, 1\"),\n    @(\"31\u00f77=4, 3\", \"33\u00f76=5, 3\"),\n    @(\"87\u00f79=9, 6\", \"95\u00f73=31, 2\"),\n    @(\"41\u00f76=6, 5\", \"32\u00f78=4, 0\"),\n    @(\"21\u00f73=7, 0\", \"12\u00f74=3, 0\"),\n    @(\"89\u00f79=9, 8\", \"61\u00f75=12, 1\"),\n    @(\"86\u00f78=10, 6\", \"90\u00f75=18, 0\"),\n    @(\"99\u00f73=33, 0\", \"58\u00f77=8, 2\"),\n    @(\"11\u00f75=2, 1\", \"93\u00f76=15, 3\"),\n    @(\"33\u00f77=4, 5\", \"16\u00f72=8, 0\"),\n    @(\"14\u00f76=2, 2\", \"34\u00f76=5, 4\"),\n    @(\"44\u00f72=22, 0\", \"38\u00f76=6, 2\"),\n    @(\"19\u00f76=3, 1\", \"72\u00f72=36, 0\"),\n    @(\"22\u00f77=3, 1\", \"51\u00f78=6, 3\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2) | Out-Null\n}\n"}
